$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "2nd Pick" text for every response row (E2:E11): Travis Hunter's
#    projected position changes from ATH to WR/CB. Also re-style that block the
#    way the author's sheet ended up (plain Arial, centered, no table border/fill)
#    instead of the banded table look the other "2nd Pick" free-for-all cells keep.
$srcFormat = $ws.Range("B2")
$pickRange = $ws.Range("E2:E11")
$srcFormat.Copy()
$pickRange.PasteSpecial(-4122)
$pickRange.HorizontalAlignment = -4108
$pickRange.VerticalAlignment = -4107
$pickRange.Value = "Travis Hunter -- Colorado Jr WR/CB"

# 2) Widen columns D, E and G so the longer entries fit.
$ws.Columns.Item(4).ColumnWidth = 23
$ws.Columns.Item(5).ColumnWidth = 23.75
$ws.Columns.Item(7).ColumnWidth = 23.75

# 3) Stray formatted cell below the table (G23), carried over from the source
#    spreadsheet -- same "table body" look as the rest of the grid.
$ws.Range("F2").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
